$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9034742712974548
$ws.Range("B1").Value = 1.630229592323303
$ws.Range("C1").Value = 4.293648719787598
$ws.Range("D1").Value = 2.787139892578125
$ws.Range("E1").Value = 0.7380626797676086
